$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bump the "Förändrad" (changed) date column for every data row (2-157) by one day.
$ws.Range("C2:C157").Value = 45184

# Row 2 hyperlink formulas gain a friendly-name second argument.
$ws.Range("S2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_SKELLEFTEA/artfynd/A 30703-2023.xlsx, "A 30703-2023""'
$ws.Range("T2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_SKELLEFTEA/kartor/A 30703-2023.png", "A 30703-2023")'
$ws.Range("V2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_SKELLEFTEA/klagomål/A 30703-2023.docx", "A 30703-2023")'
$ws.Range("W2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_SKELLEFTEA/klagomålsmail/A 30703-2023.docx", "A 30703-2023")'
$ws.Range("X2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_SKELLEFTEA/tillsyn/A 30703-2023.docx", "A 30703-2023")'
$ws.Range("Y2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_SKELLEFTEA/tillsynsmail/A 30703-2023.docx", "A 30703-2023")'
